# Update (Analyze PO & Forecast)
#
# "Forecast Comparison" sheet: MyForecast value for week W15 (row 16) changes
# "Summary" sheet: a handful of aggregate stats/labels are recalculated

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison!D16 (MyForecast, week W15): 0 -> 1 (numeric) ---
$wsForecast.Range("D16").Value = 1

# --- Summary sheet values -------------------------------------------------
# These cells store plain text in the workbook (even though some of them look
# like numbers or dates), so they are written using a leading apostrophe to
# force Excel to keep them as text instead of auto-converting to a number or
# a date serial value.

# B9  - "Total Forecast (16 Weeks)": 3 -> 5
$wsSummary.Range("B9").Formula = "'5"

# B12 - "Max Forecast": 0 -> 1
$wsSummary.Range("B12").Formula = "'1"

# B13 - "Max Forecast Week": 2025-02-16 -> 2025-05-04
$wsSummary.Range("B13").Formula = "'2025-05-04"

# B15 - "Min Forecast Week": 2025-02-23 -> 2025-02-09
$wsSummary.Range("B15").Formula = "'2025-02-09"
